# Generate Report for Handoff
# The workbook tracks localization status for e2e\a.md and e2e\b.md across
# the "Overview", "zh-cn" and "de-de" sheets. This run generates a new
# handoff for b.md: its status flips from "Handed back: in sync with en-US"
# to "Ready for handoff", new handoff xlf files are produced, and a
# version-mismatch error detail is recorded because the handback file that
# exists is stale relative to the new handoff.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$newDate = "2016-08-29 20:51:16"

$zhHandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhHandoffDate = "2016-08-29 20:51:11"

$deHandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$deHandoffDate = "2016-08-29 20:51:16"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5a18e5c2f1e74136c35d6514b2bbf5ddcbc67b1b/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/13dddd56c9ce5e3be7c862e75036b22ad3bd0b40/e2e/b.md."

# ---- Overview sheet: row 3 is b.md ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = $newDate

# Width fed to ColumnWidth is stored in the sheet XML with a +5/6 offset by
# this engine, so back that off to land on an XML-stored width of exactly 40.
$colWidthFor40 = 40 - (5 / 6)

# ---- zh-cn sheet: row 3 is b.md ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus
# Assigning the literal "False" gets auto-typed as a boolean by this COM
# layer; Copy from a cell that already holds "False" as shared-string text
# to keep the cell's stored type (t="s") consistent with the source file.
$wsZh.Range("O2").Copy($wsZh.Range("F3"))
$wsZh.Range("G3").Value = $zhHandoffFile
$wsZh.Range("H3").Value = $zhHandoffDate
$wsZh.Range("P3").Value = $errorDetail
$wsZh.Columns.Item(16).ColumnWidth = $colWidthFor40

# ---- de-de sheet: row 3 is b.md ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("O2").Copy($wsDe.Range("F3"))
$wsDe.Range("G3").Value = $deHandoffFile
$wsDe.Range("H3").Value = $deHandoffDate
$wsDe.Range("P3").Value = $errorDetail
$wsDe.Columns.Item(16).ColumnWidth = $colWidthFor40
